$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Make room for 3 extra worker rows right after the current last data row
#    (row 56). This pushes the old "firma / observaciones" footer rows
#    (61/62) down to 64/65, and Excel keeps dimension/mergeCells in sync.
# ---------------------------------------------------------------------------
$ws.Rows.Item(57).Resize(3).Insert()

# Re-apply the heavy bottom-border ("last row") look to the new last data
# row (59), using the old row 56 (which used to be the last row) as the
# format donor.
$ws.Range("B56:J56").Copy()
$ws.Range("B59:J59").PasteSpecial(-4122)

# Rows 56-58 are no longer the final row of the table, so give them the
# regular interior-row look, borrowed from row 55.
$ws.Range("B55:J55").Copy()
$ws.Range("B56:J58").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Rewrite the worker/period rows (16-59).
#    - Rows 16-56: same worker (OTTO MARIO OSPINA BERRIO / CC 73191051) as
#      before, but the 41 monthly periods now run newest-first (2009 down
#      to 1705) instead of oldest-first.
#    - Rows 57-59: three new rows for a second worker (JAIME ALFONSO
#      HOLLMAN GONZALEZ / CC 1047377698).
# ---------------------------------------------------------------------------
$rows = @(
  @(16,"73191051","OTTO MARIO OSPINA BERRIO","2009",25740,689455),
  @(17,"73191051","OTTO MARIO OSPINA BERRIO","2008",27578,689455),
  @(18,"73191051","OTTO MARIO OSPINA BERRIO","2007",27578,689455),
  @(19,"73191051","OTTO MARIO OSPINA BERRIO","2006",27578,689455),
  @(20,"73191051","OTTO MARIO OSPINA BERRIO","2005",27578,689455),
  @(21,"73191051","OTTO MARIO OSPINA BERRIO","2004",27578,689455),
  @(22,"73191051","OTTO MARIO OSPINA BERRIO","2003",27578,689455),
  @(23,"73191051","OTTO MARIO OSPINA BERRIO","2002",27578,689455),
  @(24,"73191051","OTTO MARIO OSPINA BERRIO","2001",27578,689455),
  @(25,"73191051","OTTO MARIO OSPINA BERRIO","1912",27578,689455),
  @(26,"73191051","OTTO MARIO OSPINA BERRIO","1911",27578,689455),
  @(27,"73191051","OTTO MARIO OSPINA BERRIO","1910",27578,689455),
  @(28,"73191051","OTTO MARIO OSPINA BERRIO","1909",27578,689455),
  @(29,"73191051","OTTO MARIO OSPINA BERRIO","1908",27578,689455),
  @(30,"73191051","OTTO MARIO OSPINA BERRIO","1907",27578,689455),
  @(31,"73191051","OTTO MARIO OSPINA BERRIO","1906",27578,689455),
  @(32,"73191051","OTTO MARIO OSPINA BERRIO","1905",27578,689455),
  @(33,"73191051","OTTO MARIO OSPINA BERRIO","1904",27578,689455),
  @(34,"73191051","OTTO MARIO OSPINA BERRIO","1903",27578,689455),
  @(35,"73191051","OTTO MARIO OSPINA BERRIO","1902",27578,689455),
  @(36,"73191051","OTTO MARIO OSPINA BERRIO","1901",27578,689455),
  @(37,"73191051","OTTO MARIO OSPINA BERRIO","1812",27578,689455),
  @(38,"73191051","OTTO MARIO OSPINA BERRIO","1811",27578,689455),
  @(39,"73191051","OTTO MARIO OSPINA BERRIO","1810",27578,689455),
  @(40,"73191051","OTTO MARIO OSPINA BERRIO","1809",27578,689455),
  @(41,"73191051","OTTO MARIO OSPINA BERRIO","1808",27578,689455),
  @(42,"73191051","OTTO MARIO OSPINA BERRIO","1807",27578,689455),
  @(43,"73191051","OTTO MARIO OSPINA BERRIO","1806",27578,689455),
  @(44,"73191051","OTTO MARIO OSPINA BERRIO","1805",27578,689455),
  @(45,"73191051","OTTO MARIO OSPINA BERRIO","1804",27578,689455),
  @(46,"73191051","OTTO MARIO OSPINA BERRIO","1803",27578,689455),
  @(47,"73191051","OTTO MARIO OSPINA BERRIO","1802",27578,689455),
  @(48,"73191051","OTTO MARIO OSPINA BERRIO","1801",27578,689455),
  @(49,"73191051","OTTO MARIO OSPINA BERRIO","1712",27578,689455),
  @(50,"73191051","OTTO MARIO OSPINA BERRIO","1711",27578,689455),
  @(51,"73191051","OTTO MARIO OSPINA BERRIO","1710",27578,689455),
  @(52,"73191051","OTTO MARIO OSPINA BERRIO","1709",27578,689455),
  @(53,"73191051","OTTO MARIO OSPINA BERRIO","1708",27578,689455),
  @(54,"73191051","OTTO MARIO OSPINA BERRIO","1707",27578,689455),
  @(55,"73191051","OTTO MARIO OSPINA BERRIO","1706",27578,689455),
  @(56,"73191051","OTTO MARIO OSPINA BERRIO","1705",27578,689455),
  @(57,"1047377698","JAIME ALFONSO HOLLMAN GONZALEZ","1807",31249,781242),
  @(58,"1047377698","JAIME ALFONSO HOLLMAN GONZALEZ","1806",31249,781242),
  @(59,"1047377698","JAIME ALFONSO HOLLMAN GONZALEZ","1802",27578,781242)
)

foreach ($row in $rows) {
  $r = $row[0]
  $doc = $row[1]
  $name = $row[2]
  $period = $row[3]
  $valorMora = $row[4]
  $salario = $row[5]

  $ws.Cells.Item($r, 2).Value = "CC"
  $ws.Cells.Item($r, 3).Value = $doc
  $ws.Cells.Item($r, 4).Value = $name
  $ws.Cells.Item($r, 5).Value = $period
  $ws.Cells.Item($r, 6).Value = $valorMora
  $ws.Cells.Item($r, 7).Value = $salario
}

# ---------------------------------------------------------------------------
# 3. Update the summary box at the top of the sheet.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 1218936   # VALOR MORA total
$ws.Range("C13").Value = 2         # Cant. Trabajadores
